$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F "想去人数" (number of people wanting to go)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 238
$ws1.Range("F3").Value = 1085
$ws1.Range("F5").Value = 418
$ws1.Range("F6").Value = 76
$ws1.Range("F7").Value = 546
$ws1.Range("F9").Value = 6750
$ws1.Range("F15").Value = 1084
$ws1.Range("F16").Value = 16116
$ws1.Range("F17").Value = 1579
$ws1.Range("F22").Value = 11314
$ws1.Range("F24").Value = 933
$ws1.Range("F25").Value = 4446
$ws1.Range("F26").Value = 303
$ws1.Range("F28").Value = 43
$ws1.Range("F29").Value = 37

# Sheet "全部类型" (All types) - same underlying rows, shifted due to extra rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 238
$ws4.Range("F3").Value = 1085
$ws4.Range("F5").Value = 418
$ws4.Range("F6").Value = 76
$ws4.Range("F7").Value = 546
$ws4.Range("F10").Value = 6750
$ws4.Range("F17").Value = 1084
$ws4.Range("F18").Value = 16116
$ws4.Range("F19").Value = 1579
$ws4.Range("F26").Value = 11314
$ws4.Range("F28").Value = 933
$ws4.Range("F29").Value = 4446
$ws4.Range("F30").Value = 303
$ws4.Range("F32").Value = 43
$ws4.Range("F33").Value = 37
